# Applies the NCGNT publication-certificate letter template changes:
#   1. Split the "Прошу выдать справку ..." paragraph so that the bold
#      {{student_full_name}} placeholder becomes its own paragraph; change
#      the trailing " " before it to ":" and drop the trailing "." run.
#   2. Re-touch the "К письму прилагаю ..." run so Word re-serializes it
#      without the (unneeded) xml:space="preserve".
#   3. Split each "Label: {{placeholder}}" contact line into two runs (the
#      "Label: " text run and the placeholder run) without changing the
#      visible text.

$d = $word.ActiveDocument

function Split-ParagraphAtOffset($contextText, $prefixLen) {
    # Locates the unique occurrence of $contextText in the document body and
    # breaks the paragraph that contains it into two paragraphs at the
    # offset $prefixLen characters into the match (i.e. right before the
    # placeholder run). The two halves stay as separate paragraphs.
    $full = $d.Content
    $found = $full.Find.Execute($contextText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) {
        return $false
    }
    $boundary = $full.Start + $prefixLen
    $insPoint = $d.Range($boundary, $boundary)
    $insPoint.InsertParagraphBefore()
    return $true
}

function Split-RunAtOffset($contextText, $prefixLen) {
    # Locates the unique occurrence of $contextText in the document body and
    # splits the run that contains it into two runs at the offset
    # $prefixLen characters into the match (i.e. right before the
    # placeholder). Implemented by inserting a paragraph break at that
    # position and immediately deleting the inserted paragraph mark, which
    # leaves two sibling runs in the same paragraph.
    $full = $d.Content
    $found = $full.Find.Execute($contextText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) {
        return $false
    }
    $boundary = $full.Start + $prefixLen
    $insPoint = $d.Range($boundary, $boundary)
    $insPoint.InsertParagraphBefore()
    $markRange = $d.Range($boundary, $boundary + 1)
    $markRange.Delete()
    return $true
}

# --- Change 1: "Прошу выдать справку..." paragraph -----------------------

# 1a. " " -> ":" right before the bold {{student_full_name}} run (stays
#     inside the plain descriptive run, never touching the bold run).
$d.Content.Find.Execute("докторанта ", $true, $false, $false, $false, $false, $true, 1, $false, "докторанта:", 2) | Out-Null

# 1b. Drop the trailing "." run that used to follow {{student_full_name}}.
$d.Content.Find.Execute("{{student_full_name}}.", $true, $false, $false, $false, $false, $true, 1, $false, "{{student_full_name}}", 2) | Out-Null

# 1c. Split the paragraph so {{student_full_name}} becomes its own paragraph.
$prefix1 = "для PhD докторанта:"
Split-ParagraphAtOffset "для PhD докторанта:{{student_full_name}}" $prefix1.Length | Out-Null

# --- Change 2: re-touch "К письму прилагаю..." run -----------------------
# No-op replace forces the run to be re-emitted; the serializer only keeps
# xml:space="preserve" when the text actually needs it (no leading/trailing
# whitespace here), so this drops the now-superfluous attribute.
$d.Content.Find.Execute("К письму прилагаю удостоверение личности в цифровом формате.", $true, $false, $false, $false, $false, $true, 1, $false, "К письму прилагаю удостоверение личности в цифровом формате.", 2) | Out-Null

# --- Change 3: split "Label: {{placeholder}}" contact lines --------------

$prefix2 = "Ф.И.О.: "
Split-RunAtOffset "Ф.И.О.: {{student_full_name}}" $prefix2.Length | Out-Null

$prefix3 = "Телефон: "
Split-RunAtOffset "Телефон: {{student_phone}}" $prefix3.Length | Out-Null

$prefix4 = "Электронная почта: "
Split-RunAtOffset "Электронная почта: {{student_email}}" $prefix4.Length | Out-Null
